# Apply a row permutation to columns D, L, M, N, O, P, Q, R, S, T for rows 2-8
# (weekly price data for Damasco got reshuffled between dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Mapping: new row number -> source row number (where the data used to live)
$rowMap = @{
    2 = 4
    3 = 7
    4 = 6
    5 = 3
    6 = 2
    7 = 8
    8 = 5
}

# Snapshot the original values for every affected cell before writing anything,
# since sources and destinations overlap.
$snapshot = @{}
foreach ($row in 2..8) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Now write the permuted values back.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $dstAddr = "$col$newRow"
        $ws.Range($dstAddr).Value = $snapshot[$srcAddr]
    }
}
